$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: extend series for new rows 20-29 (values 18-29) ---
# Copy style from A19 (bold/border/center) down to A20:A31
$ws.Range("A19").Copy() | Out-Null
$ws.Range("A20:A31").PasteSpecial(-4122) | Out-Null
for ($r = 20; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# --- Columns B:F content updates ---
$ws.Range("B2").Value = "NSE:AGSTRA"
$ws.Range("C2").Value = "NSE:DCM"
$ws.Range("D2").Value = "NSE:BEL"
$ws.Range("E2").ClearContents() | Out-Null
$ws.Range("F2").Value = "NSE:INDUSTOWER"
$ws.Range("B3").Value = "NSE:APOLLO"
$ws.Range("C3").Value = "NSE:MUTHOOTCAP"
$ws.Range("D3").ClearContents() | Out-Null
$ws.Range("E3").ClearContents() | Out-Null
$ws.Range("F3").Value = "NSE:NATIONALUM"
$ws.Range("B4").Value = "NSE:ASHAPURMIN"
$ws.Range("C4").ClearContents() | Out-Null
$ws.Range("D4").ClearContents() | Out-Null
$ws.Range("E4").ClearContents() | Out-Null
$ws.Range("F4").ClearContents() | Out-Null
$ws.Range("B5").Value = "NSE:ATAM"
$ws.Range("C5").ClearContents() | Out-Null
$ws.Range("D5").ClearContents() | Out-Null
$ws.Range("E5").ClearContents() | Out-Null
$ws.Range("F5").ClearContents() | Out-Null
$ws.Range("B6").Value = "NSE:CCHHL"
$ws.Range("C6").ClearContents() | Out-Null
$ws.Range("D6").ClearContents() | Out-Null
$ws.Range("E6").ClearContents() | Out-Null
$ws.Range("F6").ClearContents() | Out-Null
$ws.Range("B7").Value = "NSE:COUNCODOS"
$ws.Range("C7").ClearContents() | Out-Null
$ws.Range("D7").ClearContents() | Out-Null
$ws.Range("E7").ClearContents() | Out-Null
$ws.Range("F7").ClearContents() | Out-Null
$ws.Range("B8").Value = "NSE:CPSEETF"
$ws.Range("C8").ClearContents() | Out-Null
$ws.Range("D8").ClearContents() | Out-Null
$ws.Range("E8").ClearContents() | Out-Null
$ws.Range("F8").ClearContents() | Out-Null
$ws.Range("B9").Value = "NSE:DEN"
$ws.Range("C9").ClearContents() | Out-Null
$ws.Range("D9").ClearContents() | Out-Null
$ws.Range("E9").ClearContents() | Out-Null
$ws.Range("F9").ClearContents() | Out-Null
$ws.Range("B10").Value = "NSE:DHANUKA"
$ws.Range("C10").ClearContents() | Out-Null
$ws.Range("D10").ClearContents() | Out-Null
$ws.Range("E10").ClearContents() | Out-Null
$ws.Range("F10").ClearContents() | Out-Null
$ws.Range("B11").Value = "NSE:GOCOLORS"
$ws.Range("C11").ClearContents() | Out-Null
$ws.Range("D11").ClearContents() | Out-Null
$ws.Range("E11").ClearContents() | Out-Null
$ws.Range("F11").ClearContents() | Out-Null
$ws.Range("B12").Value = "NSE:GULFPETRO"
$ws.Range("C12").ClearContents() | Out-Null
$ws.Range("D12").ClearContents() | Out-Null
$ws.Range("E12").ClearContents() | Out-Null
$ws.Range("F12").ClearContents() | Out-Null
$ws.Range("B13").Value = "NSE:INDIAGLYCO"
$ws.Range("C13").ClearContents() | Out-Null
$ws.Range("D13").ClearContents() | Out-Null
$ws.Range("E13").ClearContents() | Out-Null
$ws.Range("F13").ClearContents() | Out-Null
$ws.Range("B14").Value = "NSE:INDUSTOWER"
$ws.Range("C14").ClearContents() | Out-Null
$ws.Range("D14").ClearContents() | Out-Null
$ws.Range("E14").ClearContents() | Out-Null
$ws.Range("F14").ClearContents() | Out-Null
$ws.Range("B15").Value = "NSE:IVC"
$ws.Range("C15").ClearContents() | Out-Null
$ws.Range("D15").ClearContents() | Out-Null
$ws.Range("E15").ClearContents() | Out-Null
$ws.Range("F15").ClearContents() | Out-Null
$ws.Range("B16").Value = "NSE:JINDRILL"
$ws.Range("C16").ClearContents() | Out-Null
$ws.Range("D16").ClearContents() | Out-Null
$ws.Range("E16").ClearContents() | Out-Null
$ws.Range("F16").ClearContents() | Out-Null
$ws.Range("B17").Value = "NSE:JISLDVREQS"
$ws.Range("C17").ClearContents() | Out-Null
$ws.Range("D17").ClearContents() | Out-Null
$ws.Range("E17").ClearContents() | Out-Null
$ws.Range("F17").ClearContents() | Out-Null
$ws.Range("B18").Value = "NSE:JISLJALEQS"
$ws.Range("C18").ClearContents() | Out-Null
$ws.Range("D18").ClearContents() | Out-Null
$ws.Range("E18").ClearContents() | Out-Null
$ws.Range("F18").ClearContents() | Out-Null
$ws.Range("B19").Value = "NSE:KIRIINDUS"
$ws.Range("C19").ClearContents() | Out-Null
$ws.Range("D19").ClearContents() | Out-Null
$ws.Range("E19").ClearContents() | Out-Null
$ws.Range("F19").ClearContents() | Out-Null
$ws.Range("B20").Value = "NSE:MAHKTECH"
$ws.Range("C20").ClearContents() | Out-Null
$ws.Range("D20").ClearContents() | Out-Null
$ws.Range("E20").ClearContents() | Out-Null
$ws.Range("F20").ClearContents() | Out-Null
$ws.Range("B21").Value = "NSE:MAZDOCK"
$ws.Range("C21").ClearContents() | Out-Null
$ws.Range("D21").ClearContents() | Out-Null
$ws.Range("E21").ClearContents() | Out-Null
$ws.Range("F21").ClearContents() | Out-Null
$ws.Range("B22").Value = "NSE:MEDICO"
$ws.Range("C22").ClearContents() | Out-Null
$ws.Range("D22").ClearContents() | Out-Null
$ws.Range("E22").ClearContents() | Out-Null
$ws.Range("F22").ClearContents() | Out-Null
$ws.Range("B23").Value = "NSE:MONARCH"
$ws.Range("C23").ClearContents() | Out-Null
$ws.Range("D23").ClearContents() | Out-Null
$ws.Range("E23").ClearContents() | Out-Null
$ws.Range("F23").ClearContents() | Out-Null
$ws.Range("B24").Value = "NSE:MRPL"
$ws.Range("C24").ClearContents() | Out-Null
$ws.Range("D24").ClearContents() | Out-Null
$ws.Range("E24").ClearContents() | Out-Null
$ws.Range("F24").ClearContents() | Out-Null
$ws.Range("B25").Value = "NSE:MTARTECH"
$ws.Range("C25").ClearContents() | Out-Null
$ws.Range("D25").ClearContents() | Out-Null
$ws.Range("E25").ClearContents() | Out-Null
$ws.Range("F25").ClearContents() | Out-Null
$ws.Range("B26").Value = "NSE:MTNL"
$ws.Range("C26").ClearContents() | Out-Null
$ws.Range("D26").ClearContents() | Out-Null
$ws.Range("E26").ClearContents() | Out-Null
$ws.Range("F26").ClearContents() | Out-Null
$ws.Range("B27").Value = "NSE:PATELENG"
$ws.Range("C27").ClearContents() | Out-Null
$ws.Range("D27").ClearContents() | Out-Null
$ws.Range("E27").ClearContents() | Out-Null
$ws.Range("F27").ClearContents() | Out-Null
$ws.Range("B28").Value = "NSE:PFS"
$ws.Range("C28").ClearContents() | Out-Null
$ws.Range("D28").ClearContents() | Out-Null
$ws.Range("E28").ClearContents() | Out-Null
$ws.Range("F28").ClearContents() | Out-Null
$ws.Range("B29").Value = "NSE:PILITA"
$ws.Range("C29").ClearContents() | Out-Null
$ws.Range("D29").ClearContents() | Out-Null
$ws.Range("E29").ClearContents() | Out-Null
$ws.Range("F29").ClearContents() | Out-Null
$ws.Range("B30").Value = "NSE:PREMEXPLN"
$ws.Range("C30").ClearContents() | Out-Null
$ws.Range("D30").ClearContents() | Out-Null
$ws.Range("E30").ClearContents() | Out-Null
$ws.Range("F30").ClearContents() | Out-Null
$ws.Range("B31").Value = "NSE:PRIVISCL"
$ws.Range("C31").ClearContents() | Out-Null
$ws.Range("D31").ClearContents() | Out-Null
$ws.Range("E31").ClearContents() | Out-Null
$ws.Range("F31").ClearContents() | Out-Null
